{"js": "// Replace the worksheet date and each \"NNN\u00d7N=\" multiplication prompt\n// with the new values from the target revision. Every old string below\n// is unique within the document, so a plain text search-and-replace on\n// context.document.body is unambiguous and order-independent.\nconst replacements = [\n  [\"2025-02-26 Wednesday\", \"2025-02-27 Thursday\"],\n  [\"954\u00d78=\", \"867\u00d73=\"],\n  [\"371\u00d74=\", \"701\u00d79=\"],\n  [\"354\u00d76=\", \"932\u00d73=\"],\n  [\"370\u00d78=\", \"213\u00d72=\"],\n  [\"912\u00d73=\", \"290\u00d73=\"],\n  [\"591\u00d78=\", \"721\u00d72=\"],\n  [\"332\u00d76=\", \"459\u00d77=\"],\n  [\"396\u00d78=\", \"909\u00d78=\"],\n  [\"868\u00d78=\", \"641\u00d72=\"],\n  [\"825\u00d77=\", \"813\u00d73=\"],\n  [\"968\u00d77=\", \"823\u00d72=\"],\n  [\"762\u00d78=\", \"684\u00d76=\"],\n  [\"307\u00d76=\", \"909\u00d76=\"],\n  [\"158\u00d74=\", \"385\u00d75=\"],\n  [\"420\u00d73=\", \"391\u00d72=\"],\n  [\"132\u00d72=\", \"254\u00d75=\"],\n  [\"247\u00d79=\", \"993\u00d78=\"],\n  [\"686\u00d78=\", \"432\u00d76=\"],\n  [\"890\u00d76=\", \"758\u00d78=\"],\n  [\"783\u00d78=\", \"551\u00d79=\"],\n  [\"586\u00d77=\", \"675\u00d74=\"],\n  [\"759\u00d72=\", \"318\u00d74=\"],\n  [\"467\u00d79=\", \"795\u00d77=\"],\n  [\"198\u00d73=\", \"162\u00d78=\"],\n  [\"548\u00d77=\", \"836\u00d75=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Could not find text to replace: \"${oldText}\"`);\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the worksheet date and each \"NNN\u00d7N=\" multiplication prompt with\n# the new values from the target revision. Every \"old\" string is unique\n# within the document, so Find/Execute with Replace=wdReplaceAll (2) is\n# unambiguous for each pair.\n$d = $word.ActiveDocument\n\n$replacements = @(\n  @(\"2025-02-26 Wednesday\", \"2025-02-27 Thursday\"),\n  @(\"954\u00d78=\", \"867\u00d73=\"),\n  @(\"371\u00d74=\", \"701\u00d79=\"),\n  @(\"354\u00d76=\", \"932\u00d73=\"),\n  @(\"370\u00d78=\", \"213\u00d72=\"),\n  @(\"912\u00d73=\", \"290\u00d73=\"),\n  @(\"591\u00d78=\", \"721\u00d72=\"),\n  @(\"332\u00d76=\", \"459\u00d77=\"),\n  @(\"396\u00d78=\", \"909\u00d78=\"),\n  @(\"868\u00d78=\", \"641\u00d72=\"),\n  @(\"825\u00d77=\", \"813\u00d73=\"),\n  @(\"968\u00d77=\", \"823\u00d72=\"),\n  @(\"762\u00d78=\", \"684\u00d76=\"),\n  @(\"307\u00d76=\", \"909\u00d76=\"),\n  @(\"158\u00d74=\", \"385\u00d75=\"),\n  @(\"420\u00d73=\", \"391\u00d72=\"),\n  @(\"132\u00d72=\", \"254\u00d75=\"),\n  @(\"247\u00d79=\", \"993\u00d78=\"),\n  @(\"686\u00d78=\", \"432\u00d76=\"),\n  @(\"890\u00d76=\", \"758\u00d78=\"),\n  @(\"783\u00d78=\", \"551\u00d79=\"),\n  @(\"586\u00d77=\", \"675\u00d74=\"),\n  @(\"759\u00d72=\", \"318\u00d74=\"),\n  @(\"467\u00d79=\", \"795\u00d77=\"),\n  @(\"198\u00d73=\", \"162\u00d78=\"),\n  @(\"548\u00d77=\", \"836\u00d75=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $find = $d.Content.Find\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    # FindText, MatchCase, MatchWholeWord, MatchWildcards, MatchSoundsLike,\n    # MatchAllWordForms, Forward, Wrap(1=wdFindContinue), Format,\n    # ReplaceWith, Replace(2=wdReplaceAll)\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
